$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Registrar horas consumidas el Día 1 para las tareas de "baja de grupo" (fila 7)
# y "baja de alumno" (fila 9).
$ws.Range("H7").Value = 2
$ws.Range("H9").Value = 2

# Actualiza la selección activa de la hoja, como quedó tras la edición.
$ws.Activate()
$ws.Range("L12").Select()
